$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''61.269.49'
$ws.Range("E2").Value = '  +0.31%  '

$ws.Range("D3").Value = '''2.932.03'
$ws.Range("E3").Value = '  +0.34%  '

$ws.Range("D4").Value = '''0.999'
$ws.Range("E4").Value = '  -0.17%  '

$ws.Range("D5").Value = '''593.58'
$ws.Range("E5").Value = '  +0.43%  '

$ws.Range("D6").Value = '''145.10'
$ws.Range("E6").Value = '  -0.12%  '

$ws.Range("E7").Value = '  -0.17%  '

$ws.Range("E8").Value = '  -0.42%  '

$ws.Range("E9").Value = '  +3.31%  '

$ws.Range("D10").Value = '''0.143'
$ws.Range("E10").Value = '  -0.90%  '

$ws.Range("D11").Value = '''0.442'
$ws.Range("E11").Value = '  -0.34%  '

$ws.Range("D12").Value = '''0.0000226'
$ws.Range("E12").Value = '  -0.64%  '

$ws.Range("D13").Value = '''33.63'
$ws.Range("E13").Value = '  -0.10%  '

$ws.Range("E14").Value = '  +0.26%  '

$ws.Range("D15").Value = '''3.415.53'
$ws.Range("E15").Value = '  +0.10%  '

$ws.Range("D16").Value = '''61.189.91'
$ws.Range("E16").Value = '  +0.14%  '

$ws.Range("E17").Value = '  +0.02%  '

$ws.Range("D18").Value = '''2.932.30'
$ws.Range("E18").Value = '  +0.23%  '

$ws.Range("D19").Value = '''433.42'
$ws.Range("E19").Value = '  +0.46%  '

$ws.Range("D20").Value = '''13.50'
$ws.Range("E20").Value = '  -0.17%  '

$ws.Range("D21").Value = '''0.680'
$ws.Range("E21").Value = '  -0.39%  '

$ws.Range("D22").Value = '''7.13'
$ws.Range("E22").Value = '  +0.57%  '

$ws.Range("D23").Value = '''82.04'
$ws.Range("E23").Value = '  +1.32%  '

$ws.Range("E24").Value = '  +1.18%  '

$ws.Range("D25").Value = '''2.21'
$ws.Range("E25").Value = '  -0.62%  '

$ws.Range("D26").Value = '''11.85'
$ws.Range("E26").Value = '  -2.51%  '

$ws.Range("E27").Value = '  -0.05%  '

$ws.Range("E28").Value = '  -4.63%  '

$ws.Range("D29").Value = '''2.61'
$ws.Range("E29").Value = '  -0.40%  '

$ws.Range("D30").Value = '''6.99'
$ws.Range("E30").Value = '  -1.89%  '

$ws.Range("D31").Value = '''0.111'
$ws.Range("E31").Value = '  +2.77%  '

$ws.Range("D32").Value = '''26.83'
$ws.Range("E32").Value = '  +1.14%  '

$ws.Range("E33").Value = '  -0.21%  '

$ws.Range("D34").Value = '''0.0₃0891'
$ws.Range("E34").Value = '  +3.19%  '

$ws.Range("E35").Value = '  +0.01%  '

$ws.Range("D36").Value = '''5.66'
$ws.Range("E36").Value = '  +0.55%  '

$ws.Range("D37").Value = '''3.00'
$ws.Range("E37").Value = '  -2.90%  '

$ws.Range("E38").Value = '  -0.20%  '

$ws.Range("D39").Value = '''0.124'
$ws.Range("E39").Value = '  -0.24%  '

$ws.Range("D40").Value = '''8.64'
$ws.Range("E40").Value = '  +0.19%  '

$ws.Range("D41").Value = '''42.55'
$ws.Range("E41").Value = '  +6.47%  '

$ws.Range("E42").Value = '  -2.25%  '

$ws.Range("D43").Value = '''0.0348'
$ws.Range("E43").Value = '  -0.49%  '

$ws.Range("B44").Value = 'Bittensor'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D44").Value = '''371.05'
$ws.Range("E44").Value = '  -2.59%  '

$ws.Range("B45").Value = 'Maker'
$ws.Range("C45").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D45").Value = '''2.701.56'
$ws.Range("E45").Value = '  -0.19%  '

$ws.Range("D46").Value = '''133.46'
$ws.Range("E46").Value = '  +3.04%  '

$ws.Range("D48").Value = '''23.90'
$ws.Range("E48").Value = '  -1.26%  '

$ws.Range("E49").Value = '  -1.08%  '

$ws.Range("D50").Value = '''2.01'
$ws.Range("E50").Value = '  -0.78%  '

$ws.Range("E51").Value = '  -0.38%  '
